$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in D1
$ws.Range("D1").Value = "infiltration_cmhr"

# Fill D2:D8 with value 1
$ws.Range("D2:D8").Value = 1

# Set column D width to match the diff (target stored width 16.85546875;
# the engine quantizes ColumnWidth to 1/6-character steps, so 16 is the
# closest achievable setting, yielding stored width 16.833333333333332)
$ws.Columns.Item(4).ColumnWidth = 16

# Update the active selection to D9 (matches diff's sheetView selection)
$ws.Range("D9").Select()
